$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = 42
$ws.Range("B44").Value = 42.10082268714905
$ws.Range("C44").Value = 2.863333333333333
$ws.Range("D44").Value = 77
$ws.Range("E44").Value = 34.36
$ws.Range("F44").Value = 4
$ws.Range("G44").Value = 29
$ws.Range("H44").Value = 37
$ws.Range("I44").Value = 3
$ws.Range("J44").Value = '[((11, 5, 0), (7, 12, 0)), ((6, 1, 0), (13, 13, 0)), ((11, 8, 0), (3, 2, 0)), ((3, 2, 0), (4, 12, 0)), ((6, 7, 0), (3, 2, 0)), ((16, 7, 0), (16, 5, 0)), ((7, 6, 0), (2, 12, 0)), ((7, 8, 0), (9, 8, 0)), ((10, 1, 0), (4, 12, 0)), ((13, 15, 0), (15, 1, 0)), ((16, 7, 0), (1, 9, 0)), ((15, 8, 0), (9, 10, 0)), ((12, 6, 0), (1, 1, 0)), ((4, 5, 0), (2, 12, 0)), ((15, 6, 0), (11, 15, 0)), ((8, 4, 0), (8, 15, 0)), ((6, 1, 0), (8, 4, 0)), ((16, 15, 0), (15, 6, 0)), ((9, 8, 0), (6, 3, 0)), ((8, 4, 0), (1, 1, 0)), ((1, 3, 0), (12, 6, 0)), ((6, 8, 0), (2, 12, 0)), ((15, 1, 0), (10, 1, 0)), ((11, 8, 0), (6, 9, 0)), ((14, 2, 0), (12, 6, 0)), ((15, 8, 0), (6, 14, 0)), ((6, 14, 0), (9, 8, 0)), ((13, 7, 0), (6, 8, 0)), ((15, 8, 0), (3, 2, 0)), ((10, 5, 0), (15, 12, 0)), ((4, 5, 0), (8, 4, 0)), ((6, 9, 0), (15, 3, 0)), ((9, 13, 0), (10, 6, 0)), ((12, 3, 0), (11, 8, 0)), ((14, 2, 0), (4, 5, 0)), ((10, 5, 0), (6, 14, 0)), ((2, 6, 0), (15, 12, 0)), ((3, 2, 0), (4, 5, 0)), ((4, 14, 0), (6, 3, 0)), ((1, 15, 0), (12, 2, 0)), ((1, 9, 0), (12, 2, 0)), ((1, 5, 0), (16, 15, 0)), ((2, 12, 0), (8, 15, 0)), ((15, 1, 0), (4, 14, 0)), ((13, 7, 0), (13, 15, 0)), ((2, 4, 0), (1, 9, 0)), ((15, 12, 0), (4, 5, 0)), ((1, 1, 0), (8, 15, 0)), ((1, 9, 0), (15, 1, 0)), ((6, 15, 0), (12, 11, 0))]'
$ws.Range("K44").Value = 1
$ws.Range("L44").Value = 1718

$ws.Range("A45").Value = 43
$ws.Range("B45").Value = 85.43693900108337
$ws.Range("C45").Value = 2.536666666666667
$ws.Range("D45").Value = 71
$ws.Range("E45").Value = 30.44
$ws.Range("F45").Value = 4
$ws.Range("G45").Value = 17
$ws.Range("H45").Value = 48
$ws.Range("I45").Value = 3
$ws.Range("J45").Value = '[((2, 12, 0), (8, 15, 0)), ((1, 15, 0), (12, 2, 0)), ((8, 4, 0), (8, 15, 0)), ((16, 7, 0), (16, 5, 0)), ((15, 1, 0), (4, 14, 0)), ((4, 5, 0), (8, 4, 0)), ((15, 8, 0), (6, 14, 0)), ((15, 1, 0), (10, 1, 0)), ((4, 5, 0), (2, 12, 0)), ((1, 9, 0), (15, 1, 0)), ((6, 1, 0), (8, 4, 0)), ((15, 8, 0), (3, 2, 0)), ((16, 15, 0), (15, 6, 0)), ((10, 5, 0), (15, 12, 0)), ((13, 7, 0), (6, 8, 0)), ((1, 1, 0), (8, 15, 0)), ((4, 14, 0), (6, 3, 0)), ((13, 15, 0), (15, 1, 0)), ((3, 2, 0), (4, 12, 0)), ((16, 7, 0), (1, 9, 0)), ((9, 13, 0), (10, 6, 0)), ((15, 8, 0), (9, 10, 0)), ((7, 8, 0), (9, 8, 0)), ((6, 9, 0), (15, 3, 0)), ((1, 9, 0), (12, 2, 0)), ((6, 1, 0), (13, 13, 0)), ((7, 6, 0), (2, 12, 0)), ((10, 1, 0), (4, 12, 0)), ((6, 15, 0), (12, 11, 0)), ((6, 8, 0), (2, 12, 0)), ((11, 8, 0), (6, 9, 0)), ((8, 4, 0), (1, 1, 0)), ((14, 2, 0), (4, 5, 0)), ((12, 3, 0), (11, 8, 0)), ((12, 6, 0), (1, 1, 0)), ((2, 6, 0), (15, 12, 0)), ((14, 2, 0), (12, 6, 0)), ((3, 2, 0), (4, 5, 0)), ((15, 6, 0), (11, 15, 0)), ((6, 14, 0), (9, 8, 0)), ((2, 4, 0), (1, 9, 0)), ((1, 5, 0), (16, 15, 0)), ((10, 5, 0), (6, 14, 0)), ((6, 7, 0), (3, 2, 0)), ((11, 8, 0), (3, 2, 0)), ((1, 3, 0), (12, 6, 0)), ((15, 12, 0), (4, 5, 0)), ((9, 8, 0), (6, 3, 0)), ((13, 7, 0), (13, 15, 0)), ((11, 5, 0), (7, 12, 0))]'
$ws.Range("K45").Value = 1
$ws.Range("L45").Value = 1522

$ws.Range("A46").Value = 44
$ws.Range("B46").Value = 50.91376399993896
$ws.Range("C46").Value = 2.576666666666667
$ws.Range("D46").Value = 74
$ws.Range("E46").Value = 30.92
$ws.Range("F46").Value = 4
$ws.Range("G46").Value = 33
$ws.Range("H46").Value = 41
$ws.Range("I46").Value = 3
$ws.Range("J46").Value = '[((3, 2, 0), (4, 12, 0)), ((1, 9, 0), (12, 2, 0)), ((1, 5, 0), (16, 15, 0)), ((13, 7, 0), (6, 8, 0)), ((2, 12, 0), (8, 15, 0)), ((16, 7, 0), (16, 5, 0)), ((12, 3, 0), (11, 8, 0)), ((7, 8, 0), (9, 8, 0)), ((8, 4, 0), (8, 15, 0)), ((3, 2, 0), (4, 5, 0)), ((15, 8, 0), (3, 2, 0)), ((15, 1, 0), (10, 1, 0)), ((1, 15, 0), (12, 2, 0)), ((6, 8, 0), (2, 12, 0)), ((9, 13, 0), (10, 6, 0)), ((2, 4, 0), (1, 9, 0)), ((7, 6, 0), (2, 12, 0)), ((14, 2, 0), (12, 6, 0)), ((6, 7, 0), (3, 2, 0)), ((9, 8, 0), (6, 3, 0)), ((4, 14, 0), (6, 3, 0)), ((1, 3, 0), (12, 6, 0)), ((16, 7, 0), (1, 9, 0)), ((2, 6, 0), (15, 12, 0)), ((15, 12, 0), (4, 5, 0)), ((16, 15, 0), (15, 6, 0)), ((13, 15, 0), (15, 1, 0)), ((11, 5, 0), (7, 12, 0)), ((6, 9, 0), (15, 3, 0)), ((11, 8, 0), (6, 9, 0)), ((1, 1, 0), (8, 15, 0)), ((15, 1, 0), (4, 14, 0)), ((1, 9, 0), (15, 1, 0)), ((15, 6, 0), (11, 15, 0)), ((10, 5, 0), (6, 14, 0)), ((4, 5, 0), (2, 12, 0)), ((6, 14, 0), (9, 8, 0)), ((11, 8, 0), (3, 2, 0)), ((6, 1, 0), (8, 4, 0)), ((8, 4, 0), (1, 1, 0)), ((13, 7, 0), (13, 15, 0)), ((15, 8, 0), (9, 10, 0)), ((10, 5, 0), (15, 12, 0)), ((15, 8, 0), (6, 14, 0)), ((14, 2, 0), (4, 5, 0)), ((6, 1, 0), (13, 13, 0)), ((4, 5, 0), (8, 4, 0)), ((10, 1, 0), (4, 12, 0)), ((6, 15, 0), (12, 11, 0)), ((12, 6, 0), (1, 1, 0))]'
$ws.Range("K46").Value = 1
$ws.Range("L46").Value = 1546

# Copy formatting (bold font + borders + alignment) from row 43 column A to new A cells
$ws.Range("A43").Copy()
$ws.Range("A44:A46").PasteSpecial(-4122)
$excel.CutCopyMode = 0
